$d = $word.ActiveDocument

function Find-Paragraph($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) "Polish BlockM class." paragraph: add an empty "_GoBack" bookmark right
#    after the last run (i.e. at the very end of the paragraph's text, just
#    before the paragraph mark).
#
#    NOTE: this runtime has a bug where creating a zero-length Range whose
#    position is exactly "paragraph.End - 1" (the position immediately
#    before the paragraph mark) resolves to the wrong location when handed
#    to Bookmarks.Add. We work around it by temporarily inserting a dummy
#    character after the target point (so the target point is no longer the
#    "last position before the paragraph mark"), adding the bookmark there,
#    and then deleting the dummy character again.
# ---------------------------------------------------------------------------
$target = Find-Paragraph "Polish BlockM class."
$insertAt = $target.Range.End - 1

$guard = $d.Range($insertAt, $insertAt)
$guard.InsertAfter("Z")

$bmRange = $d.Range($insertAt, $insertAt)
$d.Bookmarks.Add("_GoBack", $bmRange)

$dummy = $d.Range($insertAt, $insertAt + 1)
$dummy.Text = ""

# ---------------------------------------------------------------------------
# 2) "Create connected texture for block screen." paragraph: strike through.
# ---------------------------------------------------------------------------
$target = Find-Paragraph "Create connected texture for block screen."
$target.Range.Font.StrikeThrough = $true

# ---------------------------------------------------------------------------
# 3) "Try implementing muliblocks. In to my blocks" paragraph:
#    - strike through the whole paragraph
#    - fix "muliblocks" -> "multiblocks" by inserting a "t", ending up with
#      three separate runs: "mul", "t", "iblocks"
# ---------------------------------------------------------------------------
$target = Find-Paragraph "Try implementing muliblocks. In to my blocks"

# Strike through the whole paragraph (including the paragraph mark) first;
# this correctly marks both the pPr/rPr and every run's rPr.
$target.Range.Font.StrikeThrough = $true

$pStart = $target.Range.Start

# locate "muliblocks" relative to the paragraph start
$mulStart = $pStart + "Try implementing ".Length
$mulEnd = $mulStart + "mul".Length           # end of "mul"
$origWordEnd = $mulStart + "muliblocks".Length

# Insert the missing "t" between "mul" and "iblocks"
$insertPoint = $d.Range($mulEnd, $mulEnd)
$insertPoint.InsertAfter("t")

$tEnd = $mulEnd + 1
$newWordEnd = $origWordEnd + 1

# Force run splits by toggling formatting explicitly on each piece (simply
# re-setting the same True value is a no-op and the pieces get merged back
# into one run, so flip to False and back to True to force a real split).
$mulRange = $d.Range($mulStart, $mulEnd)
$mulRange.Font.StrikeThrough = $false
$mulRange.Font.StrikeThrough = $true

$tRange = $d.Range($mulEnd, $tEnd)
$tRange.Font.StrikeThrough = $false
$tRange.Font.StrikeThrough = $true

$iblocksRange = $d.Range($tEnd, $newWordEnd)
$iblocksRange.Font.StrikeThrough = $false
$iblocksRange.Font.StrikeThrough = $true

# ---------------------------------------------------------------------------
# 4) "Implement realistic hard object physics." paragraph: remove the
#    "_GoBack" bookmark that used to sit between the two runs, and merge the
#    two runs into a single run.
#
#    Note: Word only allows one bookmark per name in a document, so the
#    Bookmarks.Add call in step 1 above already *moved* the pre-existing
#    "_GoBack" bookmark away from this paragraph to its new location (it
#    does not create a duplicate). So there is nothing left to delete here;
#    we only need to merge the two runs back into one.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Implement realistic hard object physics.", $true, $false, $false, $false, $false, $true, 1, $false, "Implement realistic hard object physics.", 2)

Write-Output "done"
